# &&& №10265 от 25.03.2024 https://2eurostore.ru/
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark rows 16-21 (2020-2023 commemorative coins) as owned/exchanged (F column 0 -> 1)
$ws.Range("F16:F21").Value = 1

# Update the active selection on the sheet to D29 (bottom-right frozen pane)
$ws.Activate()
$ws.Range("D29").Select()
